$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column R (18th column). This shifts the
# previous R:W ("X (width)" .. "belly toward high y value") block one
# column to the right, becoming S:X, and leaves a blank (but styled)
# cell at R1 ready to receive the new "dz ind" header.
$ws.Columns.Item(18).Insert() | Out-Null

# New header for the inserted column.
$ws.Range("R1").Value = "dz ind"

# Fill in the "dx ind" (L), "dy ind" (O) and "dz ind" (R) helper
# columns with difference formulas, row by row. The first data row is
# entered on its own so that it is not grouped into the shared-formula
# block, matching how Excel records a formula typed into a single cell
# followed by a fill-down over the remaining rows.
$ws.Range("L2").Formula = "=K2-J2"
$ws.Range("L3:L14").Formula = "=K3-J3"

$ws.Range("O2").Formula = "=N2-M2"
$ws.Range("O3:O14").Formula = "=N3-M3"

$ws.Range("R2").Formula = "=Q2-P2"
$ws.Range("R3:R14").Formula = "=Q3-P3"

# Leave the selection where the author left it when saving.
$ws.Range("R15").Select() | Out-Null
